$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 17.306265
$ws.Range("H2").Value = 51.918795
$ws.Range("I2").Value = 0.5463168539988408
$ws.Range("J2").Value = 0.5463168539988407
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.9225163333333333
$ws.Range("N2").Value = 2.767549
$ws.Range("O2").Value = 0.1157072327300135
$ws.Range("P2").Value = 0.1157072327300135
$ws.Range("Q2").Value = 15.965312131495
$ws.Range("R2").Value = 143.687809183455
$ws.Range("S2").Value = 0.06321281136997269
$ws.Range("T2").Value = 0.06321281136997269
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 17.306265
$ws.Range("H3").Value = 51.918795
$ws.Range("I3").Value = 0.5463168539988408
$ws.Range("J3").Value = 0.5463168539988407
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.062688333333334
$ws.Range("N3").Value = 12.188065
$ws.Range("O3").Value = 0.5095654217806198
$ws.Range("P3").Value = 0.5095654217806198
$ws.Range("Q3").Value = 70.309960909075
$ws.Range("R3").Value = 632.7896481816751
$ws.Range("S3").Value = 0.2783841781337806
$ws.Range("T3").Value = 0.2783841781337805
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 17.306265
$ws.Range("H4").Value = 51.918795
$ws.Range("I4").Value = 0.5463168539988408
$ws.Range("J4").Value = 0.5463168539988407
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.987644666666667
$ws.Range("N4").Value = 8.962934000000001
$ws.Range("O4").Value = 0.3747273454893666
$ws.Range("P4").Value = 0.3747273454893666
$ws.Range("Q4").Value = 51.70497032717
$ws.Range("R4").Value = 465.3447329445301
$ws.Range("S4").Value = 0.2047198644950875
$ws.Range("T4").Value = 0.2047198644950874
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.344413333333334
$ws.Range("H5").Value = 10.03324
$ws.Range("I5").Value = 0.1055750256186672
$ws.Range("J5").Value = 0.1055750256186672
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.9225163333333333
$ws.Range("N5").Value = 2.767549
$ws.Range("O5").Value = 0.1157072327300135
$ws.Range("P5").Value = 0.1157072327300135
$ws.Range("Q5").Value = 3.085275925417778
$ws.Range("R5").Value = 27.76748332876
$ws.Range("S5").Value = 0.01221579405973626
$ws.Range("T5").Value = 0.01221579405973626
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.344413333333334
$ws.Range("H6").Value = 10.03324
$ws.Range("I6").Value = 0.1055750256186672
$ws.Range("J6").Value = 0.1055750256186672
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.062688333333334
$ws.Range("N6").Value = 12.188065
$ws.Range("O6").Value = 0.5095654217806198
$ws.Range("P6").Value = 0.5095654217806198
$ws.Range("Q6").Value = 13.58730903117778
$ws.Range("R6").Value = 122.2857812806
$ws.Range("S6").Value = 0.05379738245887588
$ws.Range("T6").Value = 0.05379738245887587
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.344413333333334
$ws.Range("H7").Value = 10.03324
$ws.Range("I7").Value = 0.1055750256186672
$ws.Range("J7").Value = 0.1055750256186672
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.987644666666667
$ws.Range("N7").Value = 8.962934000000001
$ws.Range("O7").Value = 0.3747273454893666
$ws.Range("P7").Value = 0.3747273454893666
$ws.Range("Q7").Value = 9.991918658462223
$ws.Range("R7").Value = 89.92726792616001
$ws.Range("S7").Value = 0.03956184910005503
$ws.Range("T7").Value = 0.03956184910005502
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.02739433333333
$ws.Range("H8").Value = 33.082183
$ws.Range("I8").Value = 0.3481081203824922
$ws.Range("J8").Value = 0.3481081203824921
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.9225163333333333
$ws.Range("N8").Value = 2.767549
$ws.Range("O8").Value = 0.1157072327300135
$ws.Range("P8").Value = 0.1157072327300135
$ws.Range("Q8").Value = 10.17295138660744
$ws.Range("R8").Value = 91.55656247946699
$ws.Range("S8").Value = 0.04027862730030458
$ws.Range("T8").Value = 0.04027862730030458
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.02739433333333
$ws.Range("H9").Value = 33.082183
$ws.Range("I9").Value = 0.3481081203824922
$ws.Range("J9").Value = 0.3481081203824921
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.062688333333334
$ws.Range("N9").Value = 12.188065
$ws.Range("O9").Value = 0.5095654217806198
$ws.Range("P9").Value = 0.5095654217806198
$ws.Range("Q9").Value = 44.80086630509945
$ws.Range("R9").Value = 403.2077967458951
$ws.Range("S9").Value = 0.1773838611879634
$ws.Range("T9").Value = 0.1773838611879634
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.02739433333333
$ws.Range("H10").Value = 33.082183
$ws.Range("I10").Value = 0.3481081203824922
$ws.Range("J10").Value = 0.3481081203824921
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.987644666666667
$ws.Range("N10").Value = 8.962934000000001
$ws.Range("O10").Value = 0.3747273454893666
$ws.Range("P10").Value = 0.3747273454893666
$ws.Range("Q10").Value = 32.94593586721356
$ws.Range("R10").Value = 296.513422804922
$ws.Range("S10").Value = 0.1304456318942242
$ws.Range("T10").Value = 0.1304456318942241

Write-Output "Applied 126 cell updates"